# Generate Report for Handback
# Updates the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns for the cc2f7e37... row (row 8) on both
# the zh-cn and de-de localization-status sheets, and widens the Error Detail
# column so the new message is readable.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5f5ba79b34946e54101418c8123b1154ef649/e2e/cc2f7e37-8e79-4e7b-be5c-19aefa9b34fc.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/217c0ed586f5e564c85b19ac485889343af755cd/e2e/cc2f7e37-8e79-4e7b-be5c-19aefa9b34fc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ba5f5ba79b34946e54101418c8123b1154ef649/e2e/cc2f7e37-8e79-4e7b-be5c-19aefa9b34fc.md."

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I8"), $handbackUrl, "", "", "cc2f7e37-8e79-4e7b-be5c-19aefa9b34fc.md")
$wsZhCn.Range("J8").Value = "cc2f7e37-8e79-4e7b-be5c-19aefa9b34fc.9776d262678901915e7e07ea25ee4fb5a2074835.zh-cn.xlf"
$wsZhCn.Range("K8").Value = "2016-08-20 04:46:43"
$wsZhCn.Range("P8").Value = $errorDetail

$wsZhCn.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I8"), $handbackUrl, "", "", "cc2f7e37-8e79-4e7b-be5c-19aefa9b34fc.md")
$wsDeDe.Range("J8").Value = "cc2f7e37-8e79-4e7b-be5c-19aefa9b34fc.9776d262678901915e7e07ea25ee4fb5a2074835.de-de.xlf"
$wsDeDe.Range("K8").Value = "2016-08-20 04:46:49"
$wsDeDe.Range("P8").Value = $errorDetail

$wsDeDe.Columns.Item(16).ColumnWidth = 39.14
